# Auto-generated edit script: updates cryptos price/volume columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values that must remain stored as TEXT
# (matching the source data which was written as literal strings, not numbers).
# Pre-format the data rows as Text so assigning "332.71" etc. via .Value does not
# get reinterpreted by Excel as a Number.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.248.05"
$ws.Range("E2").Value = "  -2.55%  "
$ws.Range("D3").Value = "1.935.48"
$ws.Range("E3").Value = "  -1.24%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "332.71"
$ws.Range("E5").Value = "  +1.74%  "
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "0.4728"
$ws.Range("E7").Value = "  -4.68%  "
$ws.Range("D8").Value = "0.4067"
$ws.Range("E8").Value = "  -3.09%  "
$ws.Range("D9").Value = "53.27"
$ws.Range("E9").Value = "  +0.71%  "
$ws.Range("D10").Value = "0.08472"
$ws.Range("E10").Value = "  -8.20%  "
$ws.Range("D11").Value = "1.054"
$ws.Range("E11").Value = "  -3.64%  "
$ws.Range("E12").Value = "  -1.80%  "
$ws.Range("D13").Value = "1.931.14"
$ws.Range("E13").Value = "  -1.44%  "
$ws.Range("D14").Value = "7.564"
$ws.Range("E14").Value = "  -3.45%  "
$ws.Range("D15").Value = "6.136"
$ws.Range("E15").Value = "  -4.65%  "
$ws.Range("D16").Value = "1.006"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "90.21"
$ws.Range("D18").Value = "0.00001068"
$ws.Range("E18").Value = "  -2.68%  "
$ws.Range("D19").Value = "0.06579"
$ws.Range("E19").Value = "  -1.30%  "
$ws.Range("D20").Value = "18.29"
$ws.Range("E20").Value = "  -4.46%  "
$ws.Range("D21").Value = "1.004"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").Value = "5.798"
$ws.Range("E22").Value = "  -2.16%  "
$ws.Range("D23").Value = "28.256.85"
$ws.Range("E23").Value = "  -2.63%  "
$ws.Range("D24").Value = "11.47"
$ws.Range("E24").Value = "  -4.18%  "
$ws.Range("D25").Value = "2.293"
$ws.Range("E25").Value = "  +1.54%  "
$ws.Range("D26").Value = "2.148.51"
$ws.Range("E26").Value = "  -2.38%  "
$ws.Range("D27").Value = "154.46"
$ws.Range("E27").Value = "  -0.46%  "
$ws.Range("D28").Value = "20.20"
$ws.Range("E28").Value = "  -1.58%  "
$ws.Range("D29").Value = "2.177"
$ws.Range("E29").Value = "  -2.77%  "
$ws.Range("D30").Value = "5.792"
$ws.Range("E30").Value = "  -7.53%  "
$ws.Range("D31").Value = "123.80"
$ws.Range("E31").Value = "  -1.79%  "
$ws.Range("D32").Value = "0.9868"
$ws.Range("E32").Value = "  -4.92%  "
$ws.Range("D33").Value = "0.09627"
$ws.Range("E33").Value = "  -1.90%  "
$ws.Range("D34").Value = "1.463"
$ws.Range("E34").Value = "  -3.07%  "
$ws.Range("D35").Value = "5.597"
$ws.Range("E35").Value = "  -3.64%  "
$ws.Range("D36").Value = "3.636"
$ws.Range("E36").Value = "  -1.06%  "
$ws.Range("D37").Value = "9.223"
$ws.Range("E37").Value = "  +3.19%  "
$ws.Range("D38").Value = "0.02325"
$ws.Range("E38").Value = "  -3.94%  "
$ws.Range("D39").Value = "0.06186"
$ws.Range("E39").Value = "  -2.28%  "
$ws.Range("D40").Value = "1.243"
$ws.Range("E40").Value = "  -5.20%  "
$ws.Range("D41").Value = "0.6206"
$ws.Range("E41").Value = "  -3.32%  "
$ws.Range("E42").Value = "  -1.71%  "
$ws.Range("D43").Value = "1.003"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").Value = "0.1908"
$ws.Range("E44").Value = "  -3.22%  "
$ws.Range("D45").Value = "1.316"
$ws.Range("E45").Value = "  -2.91%  "
$ws.Range("D46").Value = "0.5922"
$ws.Range("E46").Value = "  -4.21%  "
$ws.Range("D47").Value = "12.92"
$ws.Range("E47").Value = "  -2.60%  "
$ws.Range("D48").Value = "2.052"
$ws.Range("E48").Value = "  -6.18%  "
$ws.Range("D49").Value = "3.472"
$ws.Range("E49").Value = "  +0.43%  "
$ws.Range("D50").Value = "0.06791"
$ws.Range("E50").Value = "  -2.73%  "
$ws.Range("D51").Value = "110.31"
$ws.Range("E51").Value = "  -1.67%  "
